{"js": "// Remove the trailing \"Ver no Jupiter ...\" line, the \"\u00a9 2020 ...\" copyright\n// line, and the now-redundant blank paragraph that separated them from the\n// \"LOM3101: ... (Requisito)\" requirement line above, mirroring the upstream\n// Jekyll site rebuild that dropped this footer block from the page content.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two footer paragraphs that must be removed outright.\nlet jupiterIdx = -1;\nlet copyrightIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (jupiterIdx === -1 && t.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIdx = i;\n  } else if (copyrightIdx === -1 && t.indexOf(\"Powered by Jekyll and Github pages\") !== -1) {\n    copyrightIdx = i;\n  }\n}\n\nif (jupiterIdx === -1 || copyrightIdx === -1) {\n  throw new Error(\"Could not locate the footer paragraphs to remove.\");\n}\n\n// Also drop the blank paragraph immediately preceding \"Ver no Jupiter \u2026\" so\n// the document keeps exactly one blank separator line (matching the other\n// section breaks) instead of two in a row.\nlet blankIdx = -1;\nif (jupiterIdx > 0 && items[jupiterIdx - 1].text === \"\") {\n  blankIdx = jupiterIdx - 1;\n}\n\nconst toDelete = [jupiterIdx, copyrightIdx];\nif (blankIdx !== -1) toDelete.push(blankIdx);\n\n// Delete highest index first so earlier indices stay valid as we go.\ntoDelete.sort((a, b) => b - a);\nfor (const idx of toDelete) {\n  items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" line, the \"(c) 2020 ...\" copyright\n# line, and the now-redundant blank paragraph that separated them from the\n# \"LOM3101: ... (Requisito)\" requirement line above, mirroring the upstream\n# Jekyll site rebuild that dropped this footer block from the page content.\n\n$d = $word.ActiveDocument\n\n$jupiterIndex = -1\n$copyrightIndex = -1\n$count = $d.Paragraphs.Count\n\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($jupiterIndex -eq -1 -and $t -like \"*Ver no Jupiter*\") {\n        $jupiterIndex = $i\n    } elseif ($copyrightIndex -eq -1 -and $t -like \"*Powered by Jekyll and Github pages*\") {\n        $copyrightIndex = $i\n    }\n}\n\nif ($jupiterIndex -eq -1 -or $copyrightIndex -eq -1) {\n    throw \"Could not locate the footer paragraphs to remove.\"\n}\n\n# Also drop the blank paragraph immediately preceding \"Ver no Jupiter ...\" so\n# the document keeps exactly one blank separator line instead of two in a row.\n$blankIndex = -1\nif ($jupiterIndex -gt 1) {\n    $prevText = $d.Paragraphs.Item($jupiterIndex - 1).Range.Text\n    if ($prevText -eq \"`r\" -or $prevText -eq \"\") {\n        $blankIndex = $jupiterIndex - 1\n    }\n}\n\n$toDelete = @($jupiterIndex, $copyrightIndex)\nif ($blankIndex -ne -1) {\n    $toDelete += $blankIndex\n}\n\n# Delete highest index first so earlier indices stay valid as we go.\n$toDelete = $toDelete | Sort-Object -Descending\n\nforeach ($idx in $toDelete) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
